$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 30000
$ws.Range("J10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("N10").Value = -30586
# Row 17
$ws.Range("H17").Value = 861.575
$ws.Range("J17").Value = 787.6901
$ws.Range("L17").Value = 2363.0703
$ws.Range("N17").Value = -2699.0703
# Row 97
$ws.Range("H97").Value = 1500
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4500
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -5492
# Row 132
$ws.Range("H132").Value = 318643.72
$ws.Range("I132").Value = 176158.66
$ws.Range("K132").Value = 528475.98
$ws.Range("M132").Value = -525945.98
# Row 138
$ws.Range("H138").Value = 2169.5
$ws.Range("I138").Value = 730.6857
$ws.Range("J138").Value = 2944.246
$ws.Range("K138").Value = 2192.0571
$ws.Range("L138").Value = 8832.738000000001
$ws.Range("M138").Value = 2947.9429
$ws.Range("N138").Value = -19112.738
# Row 139
$ws.Range("H139").Value = 41269.168
$ws.Range("J139").Value = 41269.168
$ws.Range("L139").Value = 41269.168
$ws.Range("N139").Value = -51549.168
# Row 140
$ws.Range("H140").Value = 48869.285
$ws.Range("J140").Value = 49590
$ws.Range("L140").Value = 49590
$ws.Range("N140").Value = -59950

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5657.532
$ws.Range("I32").Value = 4379.8936
$ws.Range("J32").Value = 9660.799999999999
$ws.Range("K32").Value = 4379.8936
$ws.Range("L32").Value = 9660.799999999999
$ws.Range("M32").Value = -4092.8936
$ws.Range("N32").Value = -10234.8
# Row 63
$ws.Range("H63").Value = 9896387
$ws.Range("I63").Value = 10657186
$ws.Range("K63").Value = 10657186
$ws.Range("M63").Value = -10656500
# Row 66
$ws.Range("H66").Value = 9896387
$ws.Range("I66").Value = 10657186
$ws.Range("K66").Value = 53285930
$ws.Range("M66").Value = -53282498
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10640421
$ws.Range("I31").Value = 1189.7576
$ws.Range("J31").Value = 35718610
$ws.Range("K31").Value = 1189.7576
$ws.Range("L31").Value = 35718610
$ws.Range("M31").Value = -894.7575999999999
$ws.Range("N31").Value = -35719200
# Row 34
$ws.Range("H34").Value = 10640421
$ws.Range("I34").Value = 1189.7576
$ws.Range("J34").Value = 35718610
$ws.Range("K34").Value = 1189.7576
$ws.Range("L34").Value = 35718610
$ws.Range("M34").Value = -987.7575999999999
$ws.Range("N34").Value = -35719014
# Row 58
$ws.Range("H58").Value = 1427.8368
$ws.Range("I58").Value = 1196.7667
$ws.Range("J58").Value = 4027.375
$ws.Range("K58").Value = 1196.7667
$ws.Range("L58").Value = 4027.375
$ws.Range("M58").Value = -993.7666999999999
$ws.Range("N58").Value = -4433.375
# Row 62
$ws.Range("H62").Value = 3774.75
$ws.Range("I62").Value = 3774.75
$ws.Range("K62").Value = 3774.75
$ws.Range("M62").Value = -3150.75
# Row 63
$ws.Range("H63").Value = 49567.75
$ws.Range("J63").Value = 49567.75
$ws.Range("L63").Value = 49567.75
$ws.Range("N63").Value = -50939.75
# Row 65
$ws.Range("H65").Value = 3774.75
$ws.Range("I65").Value = 3774.75
$ws.Range("K65").Value = 18873.75
$ws.Range("M65").Value = -15753.75
# Row 66
$ws.Range("H66").Value = 49567.75
$ws.Range("J66").Value = 49567.75
$ws.Range("L66").Value = 148703.25
$ws.Range("N66").Value = -155567.25
# Row 68
$ws.Range("H68").Value = 52127.625
$ws.Range("J68").Value = 52127.625
$ws.Range("L68").Value = 52127.625
$ws.Range("N68").Value = -53625.625
# Row 71
$ws.Range("H71").Value = 52127.625
$ws.Range("J71").Value = 52127.625
$ws.Range("L71").Value = 156382.875
$ws.Range("N71").Value = -163870.875
# Row 82
$ws.Range("H82").Value = 39100
$ws.Range("J82").Value = 39100
$ws.Range("L82").Value = 39100
$ws.Range("N82").Value = -39822
# Row 85
$ws.Range("H85").Value = 39100
$ws.Range("J85").Value = 39100
$ws.Range("L85").Value = 39100
$ws.Range("N85").Value = -41596
# Row 88
$ws.Range("H88").Value = 39800
$ws.Range("J88").Value = 39800
$ws.Range("L88").Value = 39800
$ws.Range("N88").Value = -40612
# Row 91
$ws.Range("H91").Value = 39800
$ws.Range("J91").Value = 39800
$ws.Range("L91").Value = 39800
$ws.Range("N91").Value = -42608
# Row 99
$ws.Range("H99").Value = 8005005
$ws.Range("I99").Value = 16669929
$ws.Range("J99").Value = 6613.077
$ws.Range("K99").Value = 16669929
$ws.Range("L99").Value = 6613.077
$ws.Range("M99").Value = -16668431
$ws.Range("N99").Value = -9609.077000000001
# Row 126
$ws.Range("H126").Value = 8005005
$ws.Range("I126").Value = 16669929
$ws.Range("J126").Value = 6613.077
$ws.Range("K126").Value = 50009787
$ws.Range("L126").Value = 19839.231
$ws.Range("M126").Value = -50007317
$ws.Range("N126").Value = -24779.231
# Row 132
$ws.Range("H132").Value = 1331.5193
$ws.Range("I132").Value = 691.5227
$ws.Range("J132").Value = 4851.5
$ws.Range("K132").Value = 2074.5681
$ws.Range("L132").Value = 14554.5
$ws.Range("M132").Value = 455.4319
$ws.Range("N132").Value = -19614.5
# Row 134
$ws.Range("H134").Value = 1270.6666
$ws.Range("I134").Value = 561.17645
$ws.Range("J134").Value = 3463.6365
$ws.Range("K134").Value = 1683.52935
$ws.Range("L134").Value = 10390.9095
$ws.Range("M134").Value = 851.4706499999998
$ws.Range("N134").Value = -15460.9095
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 136
$ws.Range("H136").Value = 1427.8368
$ws.Range("I136").Value = 1196.7667
$ws.Range("J136").Value = 4027.375
$ws.Range("K136").Value = 3590.300099999999
$ws.Range("L136").Value = 12082.125
$ws.Range("M136").Value = -1040.300099999999
$ws.Range("N136").Value = -17182.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 524.5
$ws.Range("I92").Value = 298
$ws.Range("J92").Value = 569.8
$ws.Range("K92").Value = 894
$ws.Range("L92").Value = 1709.4
$ws.Range("M92").Value = 354
$ws.Range("N92").Value = -4205.4
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
# Row 99
$ws.Range("H99").Value = 2625
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 5250
$ws.Range("L99").Value = 10500
$ws.Range("M99").Value = -3004
$ws.Range("N99").Value = -14992
# Row 116
$ws.Range("H116").Value = 2949.5
$ws.Range("I116").Value = 899
$ws.Range("K116").Value = 2697
$ws.Range("M116").Value = 745
# Row 120
$ws.Range("H120").Value = 3998
$ws.Range("I120").Value = 3998
$ws.Range("K120").Value = 11994
$ws.Range("M120").Value = -7156
# Row 131
$ws.Range("H131").Value = 6250804.5
$ws.Range("I131").Value = 71428880
$ws.Range("J131").Value = 852.0959
$ws.Range("K131").Value = 214286640
$ws.Range("L131").Value = 2556.2877
$ws.Range("M131").Value = -214281600
$ws.Range("N131").Value = -12636.2877
# Row 137
$ws.Range("H137").Value = 2644.9678
$ws.Range("I137").Value = 665
$ws.Range("J137").Value = 3587.8096
$ws.Range("K137").Value = 1995
$ws.Range("L137").Value = 10763.4288
$ws.Range("M137").Value = 3105
$ws.Range("N137").Value = -20963.4288

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 136
$ws.Range("H136").Value = 18533.666
$ws.Range("J136").Value = 18533.666
$ws.Range("L136").Value = 55600.99800000001
$ws.Range("N136").Value = -60700.99800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 982.7222
$ws.Range("I16").Value = 865.2353000000001
$ws.Range("J16").Value = 2980
$ws.Range("K16").Value = 865.2353000000001
$ws.Range("L16").Value = 2980
$ws.Range("M16").Value = -695.2353000000001
$ws.Range("N16").Value = -3320
# Row 68
$ws.Range("H68").Value = 652.85364
$ws.Range("I68").Value = 594.175
$ws.Range("K68").Value = 594.175
$ws.Range("M68").Value = 154.825
# Row 71
$ws.Range("H71").Value = 652.85364
$ws.Range("I71").Value = 594.175
$ws.Range("K71").Value = 2970.875
$ws.Range("M71").Value = 773.125
# Row 81
$ws.Range("H81").Value = 79932.664
$ws.Range("J81").Value = 79932.664
$ws.Range("L81").Value = 79932.664
$ws.Range("N81").Value = -81928.664
# Row 84
$ws.Range("H84").Value = 79932.664
$ws.Range("J84").Value = 79932.664
$ws.Range("L84").Value = 239797.992
$ws.Range("N84").Value = -249781.992

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 138
$ws.Range("H138").Value = 43355.332
$ws.Range("J138").Value = 43355.332
$ws.Range("L138").Value = 43355.332
$ws.Range("N138").Value = -53635.332
# Row 140
$ws.Range("H140").Value = 62552.668
$ws.Range("J140").Value = 62552.668
$ws.Range("L140").Value = 62552.668
$ws.Range("N140").Value = -72912.66800000001
# Row 141
$ws.Range("H141").Value = 43110.453
$ws.Range("J141").Value = 43110.453
$ws.Range("L141").Value = 43110.453
$ws.Range("N141").Value = -53470.453
